$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 27: "Average" / "Average" / "Avrage Fwd & Rev" labels ---
$ws.Range("F27").Value = "Average"
$ws.Range("L27").Value = "Average"
$ws.Range("O27").Value = "Avrage Fwd & Rev"

# --- Row 30: per-block averages + forward/reverse combined average ---
$ws.Range("F30").Formula = "=AVERAGE(C28:E30)"
$ws.Range("L30").Formula = "=AVERAGE(I28:K30)"
$ws.Range("O30").Formula = "=AVERAGE(F30,-L30)"

# --- Row 32: "Average" label ---
$ws.Range("F32").Value = "Average"

# --- Row 33: "Average" label ---
$ws.Range("L33").Value = "Average"

# --- Row 35: per-block averages + forward/reverse combined average ---
$ws.Range("F35").Formula = "=AVERAGE(C33:E35)"
$ws.Range("L35").Formula = "=AVERAGE(I33:K35)"
$ws.Range("O35").Formula = "=AVERAGE(F35,-L35)"

# --- Row 37 (new row): overall deviation summary across rows 28:35, bold/comma ---
$ws.Range("G37").Formula = "=(MAX(C28:E35)-MIN(C28:E35))/AVERAGE(C28:E35)*100"
$ws.Range("G30").Copy()
$ws.Range("G37").PasteSpecial(-4122)
$ws.Range("M37").Formula = "=(MAX(I28:K35)-MIN(I28:K35))/AVERAGE(I28:K35)*100"
$ws.Range("M30").Copy()
$ws.Range("M37").PasteSpecial(-4122)

# --- Row 39: "Average" / "Average" / "Average no stopping and stopping" labels ---
$ws.Range("F39").Value = "Average"
$ws.Range("L39").Value = "Average"
$ws.Range("P39").Value = "Average no stopping and stopping"

# --- Row 42: per-block averages + forward/reverse combined + no-stop/stop combined ---
$ws.Range("F42").Formula = "=AVERAGE(C40:E42)"
$ws.Range("L42").Formula = "=AVERAGE(I40:K42)"
$ws.Range("O42").Formula = "=AVERAGE(F42,-L42)"
$ws.Range("P42").Formula = "=AVERAGE(O30,O42)"

# --- Row 44: "Average" labels ---
$ws.Range("F44").Value = "Average"
$ws.Range("L44").Value = "Average"

# --- Row 47: per-block averages + forward/reverse combined + no-stop/stop combined ---
$ws.Range("F47").Formula = "=AVERAGE(C45:E47)"
$ws.Range("L47").Formula = "=AVERAGE(I45:K47)"
$ws.Range("O47").Formula = "=AVERAGE(F47,-L47)"
$ws.Range("P47").Formula = "=AVERAGE(O35,O47)"

# --- Row 49 (new row): overall deviation summary across rows 40:47, bold/comma ---
$ws.Range("G49").Formula = "=(MAX(C40:E47)-MIN(C40:E47))/AVERAGE(C40:E47)*100"
$ws.Range("G42").Copy()
$ws.Range("G49").PasteSpecial(-4122)
$ws.Range("M49").Formula = "=(MAX(I40:K47)-MIN(I40:K47))/AVERAGE(I40:K47)*100"
$ws.Range("M42").Copy()
$ws.Range("M49").PasteSpecial(-4122)

# --- sheet view: scroll so row 16 is at top, and move the active selection ---
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G22").Select()
